# Germany Landesliga - fixture data correction (01-06-2024 01:16 update)
#
# Three pairs of adjacent rows had their match results/odds attached to the
# wrong fixture row (the "id"/"Div"/"Date" columns - A, C, D - were correct,
# but everything from column B onward, i.e. match id, HomeTeam, AwayTeam,
# score and all odds columns, belonged to the other row of the pair).
# This script swaps columns B and E:AD between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colNames = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(11, 12),
    @(88, 90),
    @(118, 119)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $colNames) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
